$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Build_URL")
$ws.Name = "URL"
$ws.Range("E34").Select()
